# Cryptocurrency price/volume refresh (GitHub Actions scheduled update).
# D column = "Price", E column = "Volume(1h)" percent-change text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.168.46"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "1.870.29"
$ws.Range("E3").Value = "  +2.55%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'311.69"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  -1.03%  "
$ws.Range("D8").Value = "'0.3928"
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("D9").Value = "'0.09678"
$ws.Range("E9").Value = "  -3.70%  "
$ws.Range("E10").Value = "  +2.50%  "
$ws.Range("D11").Value = "'40.86"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").Value = "'6.523"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").Value = "'20.91"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").Value = "1.865.45"
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("D15").Value = "'1.001"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "'7.410"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "'0.00001128"
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("D18").Value = "'92.79"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").Value = "'0.06585"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "'17.54"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").Value = "'6.162"
$ws.Range("E22").Value = "  +2.02%  "
$ws.Range("D23").Value = "28.232.37"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "'11.35"
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("D25").Value = "'2.274"
$ws.Range("E25").Value = "  +1.46%  "
$ws.Range("D26").Value = "'2.536"
$ws.Range("E26").Value = "  +3.21%  "
$ws.Range("D27").Value = "'21.27"
$ws.Range("E27").Value = "  +2.11%  "
$ws.Range("D28").Value = "2.081.47"
$ws.Range("E28").Value = "  +2.11%  "
$ws.Range("D29").Value = "'158.22"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "'127.46"
$ws.Range("E30").Value = "  -0.96%  "
$ws.Range("E31").Value = "  -2.85%  "
$ws.Range("D32").Value = "'1.067"
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").Value = "'3.626"
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("D35").Value = "'9.582"
$ws.Range("E35").Value = "  +4.90%  "
$ws.Range("D36").Value = "'0.06717"
$ws.Range("E36").Value = "  -3.50%  "
$ws.Range("D37").Value = "'0.02390"
$ws.Range("E37").Value = "  +1.95%  "
$ws.Range("D38").Value = "'0.2180"
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("D40").Value = "'0.6362"
$ws.Range("E40").Value = "  +1.41%  "
$ws.Range("D41").Value = "'4.964"
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("D42").Value = "'1.179"
$ws.Range("E42").Value = "  +1.90%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("E44").Value = "  +1.82%  "
$ws.Range("D45").Value = "'0.6009"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("E46").Value = "  -1.68%  "
$ws.Range("D47").Value = "'1.257"
$ws.Range("E47").Value = "  -2.23%  "
$ws.Range("D48").Value = "'124.17"
$ws.Range("E48").Value = "  -1.19%  "
$ws.Range("D49").Value = "'1.991"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("E51").Value = "  +0.75%  "
